$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 with the new workflow entry
$ws.Range("A5").Value = 45898
$ws.Range("A5").NumberFormat = "m/d/yy"

$ws.Range("B5").Value = "The Daily News Digest Bot"
$ws.Range("D5").Value = "The Daily News Digest Bot.json"
$ws.Range("C5").Value = "this workflow automates the entire process of news using  telegram scheduled timing and triggered message "


# Update the active selection like Excel would after editing near the new row
$ws.Range("C11").Select()
